$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row order/content (rows 2-8), reflecting de-duplication + added records.
$data = @(
    @("Water Ninjutsu Scroll[Ninjutsu Scroll]", "Ailment Resistance %`n5`nMagic Device only:Aggro %-10", "empty", "Sell1 SpinaProcess2 Wood"),
    @("Earth Ninjutsu Scroll[Ninjutsu Scroll]", "MaxHP %`n101`n-Handed Sword only:Fractional Barrier %10", "empty", "Sell1 SpinaProcess2 Wood"),
    @("Wind Ninjutsu Scroll[Ninjutsu Scroll]", "ASPD`n250`nKatana only:Critical Rate5", "empty", "Sell1 SpinaProcess2 Wood"),
    @("Metal Ninjutsu Scroll[Ninjutsu Scroll]", "Critical Rate5", "empty", "Sell1 SpinaProcess2 Wood"),
    @("Lightning Ninjutsu Scroll[Ninjutsu Scroll]", "Stability %`n5`nKatana only:Accuracy %10", "empty", "Sell1 SpinaProcess2 Wood"),
    @("Fire Ninjutsu Scroll[Ninjutsu Scroll]", "MATK %`n1`nStaff only:Magic Pierce %5", "empty", "Sell1 SpinaProcess2 Wood"),
    @("Dark Ninjutsu Scroll[Ninjutsu Scroll]", "Aggro %-10", "empty", "Sell1 SpinaProcess2 Wood")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Remove the now-unused rows 9 and 10 (previous 10-row table shrinks to 8).
$ws.Range("A9:E10").Delete()
